$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value (values published as text strings
# in the source data, e.g. "490" rather than numeric 490).
$updates = @{
    "C7" = "490"
    "D7" = "1245509.82"
    "C8" = "336"
    "D8" = "713041.00"
    "C9" = "860"
    "D9" = "3348728.81"
    "C11" = "17"
    "D11" = "43500.00"
    "C13" = "189"
    "D13" = "573948.18"
    "C14" = "112"
    "D14" = "274745.00"
    "C17" = "140"
    "D17" = "640217.26"
    "C18" = "205"
    "D18" = "464769.87"
    "C51" = "122"
    "D51" = "309228.33"
    "C52" = "18"
    "D52" = "62300.00"
    "C53" = "275"
    "D53" = "1042723.80"
    "C59" = "24"
    "D59" = "57000.00"
    "C61" = "46"
    "D61" = "212006.00"
    "C62" = "67"
    "D62" = "152500.00"
    "C80" = "221"
    "D80" = "613075.19"
    "C107" = "22"
    "D107" = "51909.00"
    "C109" = "84"
    "D109" = "208310.00"
    "C110" = "55"
    "D110" = "190772.01"
    "C111" = "22"
    "D111" = "69913.61"
    "C112" = "100"
    "D112" = "611906.82"
    "C114" = "7"
    "D114" = "19910.00"
    "C115" = "33"
    "D115" = "112910.00"
    "C119" = "22"
    "D119" = "113068.92"
    "C125" = "152"
    "D125" = "391879.45"
    "C126" = "550"
    "D126" = "2620098.06"
    "C128" = "18"
    "D128" = "59872.00"
    "C130" = "96"
    "D130" = "295740.68"
    "C143" = "2702"
    "D143" = "11831211.82"
    "C206" = "681"
    "D206" = "2646482.58"
    "C207" = "29"
    "D207" = "110238.00"
    "C246" = "24"
    "D246" = "68095.28"
    "C247" = "93"
    "D247" = "274838.00"
    "C248" = "169"
    "D248" = "428700.00"
    "C249" = "566"
    "D249" = "1491991.59"
    "C250" = "131"
    "D250" = "373627.11"
    "C251" = "1071"
    "D251" = "4047299.06"
    "C252" = "50"
    "D252" = "138503.19"
    "C253" = "26"
    "D253" = "64900.00"
    "C254" = "84"
    "D254" = "196500.00"
    "C255" = "191"
    "D255" = "600068.19"
    "C256" = "143"
    "D256" = "464393.00"
    "C257" = "104"
    "D257" = "275972.92"
    "C259" = "148"
    "D259" = "556875.09"
    "C260" = "229"
    "D260" = "523863.00"
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.Style = $origStyle
}
